$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Handout master: bump the cached "today" date field by one day
#    (17.05.2021 -> 18.05.2021)
# ---------------------------------------------------------------------
$hm = $p.HandoutMaster
$hf = $hm.HeadersFooters
$dt = $hf.DateAndTime
$dt.Text = "18.05.2021"

# ---------------------------------------------------------------------
# 2) Slide layout "Standardseite": swap the Font Awesome bullet glyphs
#    used by the body placeholder (idx 13 / "Textplatzhalter 4") for
#    plain-text bullets on levels 1-3.
# ---------------------------------------------------------------------
$layout = $p.SlideMaster.CustomLayouts.Item(2)
$layoutShape = $layout.Shapes.Item(4)
$layoutRange = $layoutShape.TextFrame.TextRange

$lvl1 = $layoutRange.Paragraphs(1, 1)
$lvl1Bullet = $lvl1.ParagraphFormat.Bullet
$lvl1Bullet.Font.Name = "Arial"
$lvl1Bullet.Character = 8226

$lvl2 = $layoutRange.Paragraphs(2, 1)
$lvl2Bullet = $lvl2.ParagraphFormat.Bullet
$lvl2Bullet.Font.Name = "Courier New"
$lvl2Bullet.Character = 111

$lvl3 = $layoutRange.Paragraphs(3, 1)
$lvl3Bullet = $lvl3.ParagraphFormat.Bullet
$lvl3Bullet.Font.Name = "Wingdings"
$lvl3Bullet.Character = 167

# ---------------------------------------------------------------------
# 3) Slide 4 ("JS Fehlerobjekt"): give the four bullet paragraphs an
#    explicit Arial "•" bullet (matching the layout change above).
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(4)
$shape = $slide.Shapes.Item(5)
$textRange = $shape.TextFrame.TextRange
$paraCount = $textRange.Paragraphs().Count

for ($i = 1; $i -le $paraCount; $i++) {
    $para = $textRange.Paragraphs($i, 1)
    $bullet = $para.ParagraphFormat.Bullet
    $bullet.Font.Name = "Arial"
    $bullet.Character = 8226
}
